$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers (row 1) ---
$ws.Range("A1").Value = "theta"
$ws.Range("B1").Value = "pixel"
$ws.Range("C1").Value = "D"

# --- New pixel (B) / D (C) data, rows 2-9 ---
$ws.Range("B2").Value = 8
$ws.Range("C2").Value = 163

$ws.Range("B3").Value = 12
$ws.Range("C3").Value = 143

$ws.Range("B4").Value = 13
$ws.Range("C4").Value = 123

$ws.Range("B5").Value = 18
$ws.Range("C5").Value = 103

$ws.Range("B6").Value = 26
$ws.Range("C6").Value = 83

$ws.Range("B7").Value = 35
$ws.Range("C7").Value = 63

$ws.Range("B8").Value = 50
$ws.Range("C8").Value = 43

$ws.Range("B9").Value = 95
$ws.Range("C9").Value = 23

# --- theta (A) formula, rows 2-9 ---
# A2 gets its own (non-shared) formula cell, A3:A9 become one shared group,
# matching how Excel lays these out when the range formula is applied in
# two steps.
$ws.Range("A2").Formula = "=ATAN2(C2,5)"
$ws.Range("A3:A9").Formula = "=ATAN2(C3,5)"

# --- Remove the old row 10 data (B10 leftover from the previous table) ---
$ws.Range("A10:C10").ClearContents()

# --- Remove the old SLOPE (row 12) / INTERCEPT (row 13) cells entirely;
#     they get rebuilt at rows 13/14 below. ---
$ws.Range("A12:C12").ClearContents()
$ws.Range("A13:C13").ClearContents()

# --- New SLOPE / INTERCEPT summary, rows 13-14 ---
$ws.Range("A13").Formula = "=SLOPE(A2:A9,B2:B9)"
$ws.Range("B13").Value = "Slope"

$ws.Range("A14").Formula = "=INTERCEPT(A2:A9,B2:B9)"
$ws.Range("B14").Value = "Intercept"

# --- Selection matches the saved state (active cell on A14) ---
$ws.Range("A14").Select() | Out-Null
